# Revert "Merge branch 'wrong-xlsform-col'"
#
# 1) The "survey" sheet's C1 header should read "message" instead of "label"
#    (the "choices" sheet keeps "label" in C1 - it was not touched by the
#    wrong merge).
# 2) The previously-active sheet/selection state ("survey", cell C2 / "choices",
#    cell C1) reverts to ("survey" selection A3, "choices" active + selection A4).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("choices")

# Content fix: survey!C1 "label" -> "message"
$ws1.Range("C1").Value = "message"

# Restore per-sheet selection state.
$ws1.Activate() | Out-Null
$ws1.Range("A3").Select() | Out-Null

# "choices" ends up the active/selected sheet, with A4 selected.
$ws2.Activate() | Out-Null
$ws2.Range("A4").Select() | Out-Null
